$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 ("Recepção solicita dados do problema"): reclassify from "Relativo*" (G)
# to "Previsível*" (E) - move the "x(3)" marker.
$ws.Range("E6").Value = $ws.Range("G6").Value()
$ws.Range("G6").ClearContents() | Out-Null

# Row 17 ("Cliente informa data e horário"): reclassify from "Não Previsível" (F)
# to "Previsível*" (E) - move the "x(12)" marker.
$ws.Range("E17").Value = $ws.Range("F17").Value()
$ws.Range("F17").ClearContents() | Out-Null

# Row 18 ("Recepção envia agenda ao técnico"): reclassify from "Relativo*" (G)
# to "Previsível*" (E) - move the "x(15)" marker.
$ws.Range("E18").Value = $ws.Range("G18").Value()
$ws.Range("G18").ClearContents() | Out-Null

# Row 28 ("Técnico envia relatório da instalação"): reclassify from "Relativo*" (G)
# to "Previsível*" (E) - move the "x(25)" marker.
$ws.Range("E28").Value = $ws.Range("G28").Value()
$ws.Range("G28").ClearContents() | Out-Null

# Row 32 ("Técnico faz a manutenção do produto"): reclassify from "Relativo*" (G)
# to "Previsível*" (E) - move the "x(29)" marker.
$ws.Range("E32").Value = $ws.Range("G32").Value()
$ws.Range("G32").ClearContents() | Out-Null

# Update the saved cursor/selection position to match the author's last position.
[void]$ws.Range("G39").Select()
